$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new "2021" column (O) mirroring the formatting of the existing
# "2020" column (N): header year in row 4, value in row 5.
$ws.Range("N4").Copy($ws.Range("O4"))
$ws.Range("O4").Value = 2021

$ws.Range("N5").Copy($ws.Range("O5"))
$ws.Range("O5").Value = 515

# Update the view: scroll back to the top-left (drop the previous
# topLeftCell="E1") and move the selection to P12.
$ws.Range("P12").Select()
